# Append a new time-tracking entry (row 9) to the "Tabelle1" sheet, mirroring
# rows 2-8: a date in column B, start/end times in C/D, a duration formula in
# E (D-C, shared with the existing E3:E8 group) and the activity label in F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: 2025-12-01, 14:30 - 17:15, activity "NanoGPT"
$ws.Range("B9").Value = 45992
$ws.Range("C9").Value = 0.60416666666666663
$ws.Range("D9").Value = 0.71875
$ws.Range("F9").Value = "NanoGPT"

# Duration formula, entered as part of a multi-cell range so the engine
# records it as a shared formula (like the existing E3:E8 block).
$ws.Range("E4:E9").Formula = "=D4-C4"

# Copy number formatting from the cell above so E9 matches the time format
# used throughout column E (this also applies the same style index, rather
# than minting a new, near-duplicate number format).
$ws.Range("E8").Copy()
$ws.Range("E9").PasteSpecial(-4122)
